$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 entirely (table shrinks from 5 data rows to 4)
$ws.Rows.Item(6).Delete()

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
$ws.Cells.Item(2, 1).Value = 0
Set-TextCell 2 2 "2025-03-19"
$ws.Cells.Item(2, 3).Value = 200
$ws.Cells.Item(2, 4).Value = "ZARAPLAST DA AMAZONIA LTDA"
Set-TextCell 2 5 "000098"
$ws.Cells.Item(2, 6).Value = "PANO DE CHAO FLANELADO C REFORCADO ITATEX 42x62CM"
$ws.Cells.Item(2, 7).Value = 939
$ws.Cells.Item(2, 8).Value = $false

# Row 3
$ws.Cells.Item(3, 1).Value = 2
Set-TextCell 3 2 "2025-03-19"
$ws.Cells.Item(3, 3).Value = 25
$ws.Cells.Item(3, 4).Value = "JURUA ESTALEIROS E NAVEGACAO LTDA"
Set-TextCell 3 5 "000276"
$ws.Cells.Item(3, 6).Value = "INSETICIDA SBP AEROSSOL 380ML"
$ws.Cells.Item(3, 7).Value = -2
$ws.Cells.Item(3, 8).Value = $false

# Row 4
$ws.Cells.Item(4, 1).Value = 3
Set-TextCell 4 2 "2025-03-19"
$ws.Cells.Item(4, 3).Value = 30
$ws.Cells.Item(4, 4).Value = "JURUA ESTALEIROS E NAVEGACAO LTDA"
Set-TextCell 4 5 "000146"
$ws.Cells.Item(4, 6).Value = "DESINFETANTE BRINORT 2L LAVANDA"
$ws.Cells.Item(4, 7).Value = 13
$ws.Cells.Item(4, 8).Value = $false

# Row 5
$ws.Cells.Item(5, 1).Value = 1
Set-TextCell 5 2 "2025-03-24"
$ws.Cells.Item(5, 3).Value = 150
$ws.Cells.Item(5, 4).Value = "JURUA ESTALEIROS E NAVEGACAO LTDA"
Set-TextCell 5 5 "000088"
$ws.Cells.Item(5, 6).Value = "VASSOURA PIACAVA 20 FUROS"
$ws.Cells.Item(5, 7).Value = -50
$ws.Cells.Item(5, 8).Value = $false
